# Update cryptocurrency price/volume data in the worksheet.
# Mirrors the upstream GitHub Actions data-refresh commit.
#
# Note: values like '''207.14' are PowerShell single-quoted literals that
# evaluate to the string  '207.14  (a leading apostrophe + the digits).
# That leading apostrophe is the classic Excel "force text" prefix, so
# numeric-looking price strings (e.g. "207.14", "1.00", "0.999") are
# stored verbatim as text instead of being coerced into real numbers
# (which would silently drop trailing/duplicate zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '80.944.80'
$ws.Range('E2').Value = '  +2.70%  '
$ws.Range('D3').Value = '3.140.15'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''207.14'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = '''617.19'
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('D7').Value = '''0.279'
$ws.Range('E7').Value = '  +23.54%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '''0.576'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '3.137.54'
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('D11').Value = '''0.571'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('E12').Value = '  +13.35%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '''5.25'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('D15').Value = '3.710.65'
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D16').Value = '''31.11'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').Value = '80.870.94'
$ws.Range('E17').Value = '  +2.72%  '
$ws.Range('D18').Value = '3.124.12'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('D19').Value = '''3.15'
$ws.Range('E19').Value = '  +11.85%  '
$ws.Range('D20').Value = '''13.75'
$ws.Range('E20').Value = '  -4.60%  '
$ws.Range('D21').Value = '''427.42'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '''8.88'
$ws.Range('E22').Value = '  -4.62%  '
$ws.Range('D23').Value = '''5.04'
$ws.Range('E23').Value = '  +2.51%  '
$ws.Range('D24').Value = '''7.15'
$ws.Range('E24').Value = '  +4.97%  '
$ws.Range('D25').Value = '''5.12'
$ws.Range('E25').Value = '  +8.65%  '
$ws.Range('D26').Value = '3.294.08'
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('D27').Value = '''75.30'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').Value = '''10.72'
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  +6.27%  '
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '''8.85'
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('D33').Value = '''556.28'
$ws.Range('E33').Value = '  +9.44%  '
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('D35').Value = '''0.142'
$ws.Range('E35').Value = '  +14.41%  '
$ws.Range('E36').Value = '  +12.10%  '
$ws.Range('D37').Value = '''1.96'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = '''22.52'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').Value = '''5.88'
$ws.Range('E41').Value = '  +9.24%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = '''20.69'
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('D43').Value = '''3.00'
$ws.Range('E43').Value = '  +21.46%  '
$ws.Range('D44').Value = '''2.00'
$ws.Range('E44').Value = '  +12.76%  '
$ws.Range('E45').Value = '  -2.81%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '''185.32'
$ws.Range('E47').Value = '  -3.10%  '
$ws.Range('D48').Value = '''44.65'
$ws.Range('E48').Value = '  +4.97%  '
$ws.Range('D49').Value = '''1.31'
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('D50').Value = '''0.760'
$ws.Range('E50').Value = '  -5.34%  '
$ws.Range('D51').Value = '''25.32'
$ws.Range('E51').Value = '  +2.09%  '
